# Updated default owning group in drools
#
# The "Sheet1" RuleTable has several "Default group" assignment rules whose
# ACTION column (G) hardcodes a literal "owning group" participant LDAP id.
# The old id pointed at the APPDEV tenant ("000.ARKCASE_*@APPDEV.ARMEDIA.COM");
# it is replaced with the production id ("ARKCASE_*@ARMEDIA.COM").
#
#   G23 - Complaint – Default group
#   G24 - Case File – Default group
#   G30 - DocumentRepository – Default group
#   G34 - Organization – Default group
#   G37 - Person – Default group

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldSupervisor = "owning group, 000.ARKCASE_SUPERVISOR@APPDEV.ARMEDIA.COM"
$newSupervisor = "owning group, ARKCASE_SUPERVISOR@ARMEDIA.COM"

$oldEntityAdmin = "owning group, 000.ARKCASE_ENTITY_ADMINISTRATOR@APPDEV.ARMEDIA.COM"
$newEntityAdmin = "owning group, ARKCASE_ENTITY_ADMINISTRATOR@ARMEDIA.COM"

$supervisorCells = @("G23", "G24", "G30")
foreach ($addr in $supervisorCells) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq $oldSupervisor) {
        $cell.Value = $newSupervisor
    }
}

$entityAdminCells = @("G34", "G37")
foreach ($addr in $entityAdminCells) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq $oldEntityAdmin) {
        $cell.Value = $newEntityAdmin
    }
}
